$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# 1) First paragraph: "- A description of your implementation of
#    BitTorrent (1 page)" -- the three runs (plain text / spell-checked
#    "BitTorrent" / plain text) collapse into a single run with no
#    proofErr spell-check markers.
# -----------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "- A description of your implementation of BitTorrent (1 page)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "- A description of your implementation of BitTorrent (1 page)", 2)

# -----------------------------------------------------------------------
# 2) Second paragraph: substantially revise the implementation
#    description wording (more accurate description of the HTTP GET to
#    the tracker, Peer objects, SocketChannel/ServerSocketChannel/
#    Selector handling, and message read/write logic), while keeping
#    the existing "peerID", "becoded" and "unchoke" spell-check runs,
#    and adding new spell-check runs around "ip", "SocketChannels",
#    "ServerSocketChannel" and "SocketChannel". We rebuild the whole
#    paragraph via InsertXML so the run / proofErr structure exactly
#    matches the revision.
# -----------------------------------------------------------------------
$p2 = $d.Paragraphs(2).Range
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00F47351" w:rsidRDefault="00F47351" w:rsidP="00F47351">
      <w:r>
        <w:t xml:space="preserve">Bit-tortoise starts off by generating our </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>peerID</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>. Then</w:t>
      </w:r>
      <w:r>
        <w:t>,</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> it </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">parses </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">the .torrent file. From the information in the torrent file we create a list of pieces and blocks inside the pieces. The purpose of the list is to organize what </w:t>
      </w:r>
      <w:r>
        <w:t>blocks we wi</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">ll </w:t>
      </w:r>
      <w:r>
        <w:t>later request and receive, to keep track of how much of a file we have gotten</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">. The next step is to </w:t>
      </w:r>
      <w:r>
        <w:t>send an HTTP GET request</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> the tracker </w:t>
      </w:r>
      <w:r>
        <w:t>to obtain a list of peers</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">. Once </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">we </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">get the response from the tracker, we parse the </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>becoded</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> data</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">, and form a list of Peer objects representing the peers </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>ip</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> addresses and ports</w:t>
      </w:r>
      <w:r>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> From the list of peers we </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">open connections (using </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>SocketChannels</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">) </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">and add them </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">and a </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>ServerSocketChannel</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">to a </w:t>
      </w:r>
      <w:r>
        <w:t>S</w:t>
      </w:r>
      <w:r>
        <w:t>elector</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">, keeping a Mapping between the </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>SocketChannel</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> and its respective peer</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">. We loop waiting for the selector to </w:t>
      </w:r>
      <w:r>
        <w:t>tell us when we have received data or are able to send data over one of the channels</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">. When a </w:t>
      </w:r>
      <w:r>
        <w:t>C</w:t>
      </w:r>
      <w:r>
        <w:t>hannel is ready</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> to be read</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">, we </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">read and </w:t>
      </w:r>
      <w:r>
        <w:t>parse the message</w:t>
      </w:r>
      <w:r>
        <w:t>, changing state in the peer object</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">. </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> When a Channel is ready to be written to, we check the state of the peer object, and send a message if necessary.  </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">Whenever we receive an </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>unchoke</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> message we request a random available piece the peer has. </w:t>
      </w:r>
    </w:p>

'@
$p2.InsertXML($xml)

Write-Host $d.Paragraphs(1).Range.Text
Write-Host $d.Paragraphs(2).Range.Text
